$d = $word.ActiveDocument

$replacements = @(
    @{old="41÷9="; new="32÷3="},
    @{old="16÷3="; new="56÷6="},
    @{old="45÷4="; new="60÷6="},
    @{old="69÷8="; new="24÷7="},
    @{old="69÷4="; new="64÷5="},
    @{old="17÷9="; new="64÷6="},
    @{old="86÷8="; new="27÷8="},
    @{old="19÷2="; new="43÷8="},
    @{old="82÷3="; new="20÷5="},
    @{old="37÷4="; new="58÷5="},
    @{old="72÷5="; new="96÷9="},
    @{old="52÷7="; new="99÷6="},
    @{old="97÷9="; new="56÷3="},
    @{old="13÷2="; new="92÷2="},
    @{old="50÷9="; new="45÷2="},
    @{old="21÷2="; new="41÷5="},
    @{old="87÷4="; new="36÷8="},
    @{old="56÷5="; new="96÷7="},
    @{old="12÷2="; new="46÷4="},
    @{old="47÷7="; new="33÷6="},
    @{old="75÷3="; new="89÷5="},
    @{old="34÷4="; new="48÷3="},
    @{old="95÷7="; new="87÷9="},
    @{old="55÷7="; new="45÷7="},
    @{old="30÷4="; new="95÷5="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
